# Update the "Waypart" (G column) raw measurements on Sheet1.
# Everything else (F8:F12/H8:H12/I8:I12 etc.) is driven by existing
# formulas and will recalculate automatically once the inputs change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G3").Value = 1598461
$ws.Range("G4").Value = 22284208
$ws.Range("G5").Value = 22296990
$ws.Range("G6").Value = 6454594

# Match the author's last cell selection on the sheet.
$ws.Range("N17").Select()
